# Update the "spike" worksheet's threshold column headers to reflect the
# long-dataframe naming used by the qaqcmar_thresholds dashboard.
$wb = $excel.ActiveWorkbook

$wsSpike = $wb.Worksheets.Item("spike")
$wsSpike.Range("B1").Value = "spike_high"
$wsSpike.Range("C1").Value = "spike_low"

# Restore the active selections on each sheet as left by the author.
$wsClimatology = $wb.Worksheets.Item("climatology")
$wsClimatology.Range("B2:B17").Select()

$wsGrossrange = $wb.Worksheets.Item("grossrange")
$wsGrossrange.Range("C2:F9").Select()

$wsSpike.Range("C2").Select()

$wsSpike.Activate()
